$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "denunciante_*" columns (C:G) entirely; the remaining
# "denuncia_*" columns (previously H:X) shift left to fill C:S.
$ws.Range("C1:G2").Delete(-4159) | Out-Null
